# Auto-generated script to apply scheduled market-price refresh to profit sheets.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H-N) per leve row
# based on the latest Universalis market data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6060652.5
$ws.Range("I11").Value = 6060652.5
$ws.Range("K11").Value = 6060652.5
$ws.Range("M11").Value = -6060512.5
$ws.Range("H64").Value = 3080.4
$ws.Range("I64").Value = 3035.5
$ws.Range("J64").Value = 3260
$ws.Range("K64").Value = 3035.5
$ws.Range("L64").Value = 3260
$ws.Range("M64").Value = -2787.5
$ws.Range("N64").Value = -3756
$ws.Range("H67").Value = 3080.4
$ws.Range("I67").Value = 3035.5
$ws.Range("J67").Value = 3260
$ws.Range("K67").Value = 3035.5
$ws.Range("L67").Value = 3260
$ws.Range("M67").Value = -2177.5
$ws.Range("N67").Value = -4976
$ws.Range("H86").Value = 2573.6365
$ws.Range("I86").Value = 2625
$ws.Range("J86").Value = 2512
$ws.Range("K86").Value = 2625
$ws.Range("L86").Value = 2512
$ws.Range("M86").Value = -1502
$ws.Range("N86").Value = -4758
$ws.Range("H89").Value = 2573.6365
$ws.Range("I89").Value = 2625
$ws.Range("J89").Value = 2512
$ws.Range("K89").Value = 13125
$ws.Range("L89").Value = 12560
$ws.Range("M89").Value = -7509
$ws.Range("N89").Value = -23792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8200.429
$ws.Range("J37").Value = 8200.429
$ws.Range("L37").Value = 8200.429
$ws.Range("N37").Value = -8746.429
$ws.Range("H63").Value = 3899.1667
$ws.Range("I63").Value = 2484.5386
$ws.Range("J63").Value = 7577.2
$ws.Range("K63").Value = 2484.5386
$ws.Range("L63").Value = 7577.2
$ws.Range("M63").Value = -1798.5386
$ws.Range("N63").Value = -8949.200000000001
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H66").Value = 3899.1667
$ws.Range("I66").Value = 2484.5386
$ws.Range("J66").Value = 7577.2
$ws.Range("K66").Value = 12422.693
$ws.Range("L66").Value = 37886
$ws.Range("M66").Value = -8990.692999999999
$ws.Range("N66").Value = -44750
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H74").Value = 972.7941
$ws.Range("I74").Value = 556.25
$ws.Range("J74").Value = 2916.6667
$ws.Range("K74").Value = 556.25
$ws.Range("L74").Value = 2916.6667
$ws.Range("M74").Value = 317.75
$ws.Range("N74").Value = -4664.6667
$ws.Range("H77").Value = 972.7941
$ws.Range("I77").Value = 556.25
$ws.Range("J77").Value = 2916.6667
$ws.Range("K77").Value = 2781.25
$ws.Range("L77").Value = 14583.3335
$ws.Range("M77").Value = 1586.75
$ws.Range("N77").Value = -23319.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2767.697
$ws.Range("I20").Value = 2050.8948
$ws.Range("J20").Value = 3740.5
$ws.Range("K20").Value = 2050.8948
$ws.Range("L20").Value = 3740.5
$ws.Range("M20").Value = -1803.8948
$ws.Range("N20").Value = -4234.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 875.2381
$ws.Range("I94").Value = 798.5
$ws.Range("J94").Value = 1335.6666
$ws.Range("K94").Value = 798.5
$ws.Range("L94").Value = 1335.6666
$ws.Range("M94").Value = -347.5
$ws.Range("N94").Value = -2237.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 70998.75
$ws.Range("I86").Value = 14198.5
$ws.Range("J86").Value = 355000
$ws.Range("K86").Value = 14198.5
$ws.Range("L86").Value = 355000
$ws.Range("M86").Value = -13075.5
$ws.Range("N86").Value = -357246
$ws.Range("H89").Value = 70998.75
$ws.Range("I89").Value = 14198.5
$ws.Range("J89").Value = 355000
$ws.Range("K89").Value = 70992.5
$ws.Range("L89").Value = 1775000
$ws.Range("M89").Value = -65376.5
$ws.Range("N89").Value = -1786232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 5671.2
$ws.Range("J43").Value = 5999.1304
$ws.Range("L43").Value = 17997.3912
$ws.Range("N43").Value = -18225.3912
$ws.Range("H107").Value = 459.5
$ws.Range("I107").Value = 286.66666
$ws.Range("J107").Value = 533.5714
$ws.Range("K107").Value = 859.9999799999999
$ws.Range("L107").Value = 1600.7142
$ws.Range("M107").Value = 1060.00002
$ws.Range("N107").Value = -5440.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 500000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H70").Value = 8935809
$ws.Range("I70").Value = 12235114
$ws.Range("J70").Value = 8278
$ws.Range("K70").Value = 12235114
$ws.Range("L70").Value = 8278
$ws.Range("M70").Value = -12234844
$ws.Range("N70").Value = -8818
$ws.Range("H73").Value = 8935809
$ws.Range("I73").Value = 12235114
$ws.Range("J73").Value = 8278
$ws.Range("K73").Value = 12235114
$ws.Range("L73").Value = 8278
$ws.Range("M73").Value = -12234178
$ws.Range("N73").Value = -10150
$ws.Range("H97").Value = 1422.3158
$ws.Range("I97").Value = 712
$ws.Range("K97").Value = 712
$ws.Range("M97").Value = -216
$ws.Range("H122").Value = 2865.75
$ws.Range("I122").Value = 3261
$ws.Range("J122").Value = 1680
$ws.Range("K122").Value = 9783
$ws.Range("L122").Value = 5040
$ws.Range("M122").Value = -7333
$ws.Range("N122").Value = -9940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2618.6365
$ws.Range("I40").Value = 2497.2
$ws.Range("J40").Value = 2878.8572
$ws.Range("K40").Value = 2497.2
$ws.Range("L40").Value = 2878.8572
$ws.Range("M40").Value = -2361.2
$ws.Range("N40").Value = -3150.8572
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
$ws.Range("H100").Value = 2593.8708
$ws.Range("J100").Value = 3356.5
$ws.Range("L100").Value = 3356.5
$ws.Range("N100").Value = -4438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1712.1538
$ws.Range("I122").Value = 1635.3
$ws.Range("J122").Value = 1968.3334
$ws.Range("K122").Value = 4905.9
$ws.Range("L122").Value = 5905.0002
$ws.Range("M122").Value = -2455.9
$ws.Range("N122").Value = -10805.0002
